$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("თელავი")

# Update the social package recipients row (row 4) for years 2015-2021 (columns E-K)
$ws.Range("E4").Value = 2299
$ws.Range("F4").Value = 2304
$ws.Range("G4").Value = 2295
$ws.Range("H4").Value = 2270
$ws.Range("I4").Value = 2336
$ws.Range("J4").Value = 2387
$ws.Range("K4").Value = 2447
